# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet named "2022-Q1" right before the "总计" (Total)
#    sheet, and fill it with the per-fund holding details for 2022-Q1.
# 2. Insert a new leading row into the "总计" sheet summarising 2022-Q1
#    (date/count/market value), shifting the previously existing rows down
#    and renumbering the index column.
# ---------------------------------------------------------------------------

function SetTextPlain($cell, $val) {
    # Use for text that can never be mistaken for a number (e.g. Chinese
    # names) - keeps whatever style/format the cell already has.
    $cell.Value = $val
}

function SetTextForced($cell, $val) {
    # Use for text that looks numeric (fund codes, decimal-looking
    # percentages, etc.) so Excel does not silently convert it into a
    # real number (which would drop leading zeros / change precision).
    # The leading apostrophe forces text interpretation; resetting the
    # style back to "Normal" clears the quote-prefix flag that gets
    # attached to the style so the cell ends up with the default style,
    # matching a plain un-styled data cell.
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

function SetIndexCell($ws, $templateAddr, $cell, $val) {
    # Column A "index" cells use the bold/centered/bordered style that is
    # already present elsewhere in the workbook. Copy that exact style in
    # (copy also brings a value which we immediately overwrite) so the
    # cell ends up byte-identical in formatting to its siblings.
    $ws.Range($templateAddr).Copy($cell)
    $cell.Value = $val
}

function SetDataRow($ws, $r, $idx, $code, $name, $scale, $pos, $pct, $mktval, $rank) {
    SetIndexCell $ws "A2" $ws.Cells.Item($r,1) $idx
    SetTextForced $ws.Cells.Item($r,2) $code
    SetTextForced $ws.Cells.Item($r,3) $name
    SetTextForced $ws.Cells.Item($r,4) $scale
    SetTextForced $ws.Cells.Item($r,5) $pos
    SetTextForced $ws.Cells.Item($r,6) $pct
    SetTextForced $ws.Cells.Item($r,7) $mktval
    $ws.Cells.Item($r,8).Value = $rank
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert the new "2022-Q1" worksheet immediately before "总计"
# ---------------------------------------------------------------------------
$lastQuarterSheet = $wb.Worksheets.Item(5)       # "2021-Q4"

$newSheet = $wb.Worksheets.Add($null, $lastQuarterSheet)
$newSheet.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet reference only *after* the new sheet has
# been inserted - worksheet references resolved by position before an
# Add()/move can end up pointing at a different sheet once the
# collection shifts, so always re-resolve by position afterwards.
$totalSheet = $wb.Worksheets.Item(7)              # "总计" (after insert)

# Copy the header row style (B1:H1, bold/centered/bordered) from an
# existing quarter sheet, then copy the column-A template style too.
$lastQuarterSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$lastQuarterSheet.Range("A2").Copy($newSheet.Range("A2"))

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    SetTextPlain $newSheet.Cells.Item(1, 2 + $i) $headers[$i]
}

SetDataRow $newSheet 2 0 "012412" "汇泉策略优选混合型证券投资基金"       "23.52" "70.35" "2.31" "0.5433" 7
SetDataRow $newSheet 3 1 "080012" "长盛电子信息产业混合"                 "6.40"  "87.51" "3.79" "0.2426" 4
SetDataRow $newSheet 4 2 "001120" "东方睿鑫热点挖掘灵活配置混合A"         "1.74"  "78.43" "4.00" "0.0696" 8
SetDataRow $newSheet 5 3 "001121" "东方睿鑫热点挖掘灵活配置混合C"         "1.10"  "78.43" "4.00" "0.0440" 8
SetDataRow $newSheet 6 4 "004205" "东方支柱产业灵活配置混合"             "0.96"  "84.14" "4.12" "0.0396" 9
SetDataRow $newSheet 7 5 "003456" "信达澳银新目标灵活配置混合"           "1.05"  "86.04" "1.98" "0.0208" 6
SetDataRow $newSheet 8 6 "000549" "华安大国新经济股票"                   "0.91"  "94.60" "1.81" "0.0165" 8
SetDataRow $newSheet 9 7 "710002" "富安达策略精选混合"                   "0.63"  "63.45" "1.53" "0.0096" 8

# ---------------------------------------------------------------------------
# Step 2: prepend a 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$a2 = $totalSheet.Cells.Item(2, 1)
SetIndexCell $totalSheet "A3" $a2 0

$b2 = $totalSheet.Cells.Item(2, 2)
SetTextForced $b2 "2022-Q1"

$c2 = $totalSheet.Cells.Item(2, 3)
$c2.Value = 8
$c2.Style = "Normal"

$d2 = $totalSheet.Cells.Item(2, 4)
$d2.Value = 0.99
$d2.Style = "Normal"

# Renumber the index column for the rows that shifted down (were 0..4,
# now need to read 1..5).
for ($i = 3; $i -le 7; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}
